# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-09-12 (serial 45181) to 2023-09-13 (serial 45182).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows run from row 2 to row 533 (row 1 is the header row).
$rng = $ws.Range("C2:C533")
$rng.Value = 45182
